$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Formula = $text
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue "D2" '="26.865.44"'
Set-TextValue "D3" '="1.808.79"'
Set-TextValue "E3" '="  -0.88%  "'
Set-TextValue "E4" '="  +0.09%  "'
Set-TextValue "D5" '="309.59"'
Set-TextValue "E5" '="  -1.38%  "'
Set-TextValue "D7" '="0.4641"'
Set-TextValue "E7" '="  -0.66%  "'
Set-TextValue "D8" '="0.3691"'
Set-TextValue "E8" '="  -2.46%  "'
Set-TextValue "D9" '="0.07348"'
Set-TextValue "E9" '="  -1.21%  "'
Set-TextValue "D10" '="0.8711"'
Set-TextValue "E10" '="  -0.48%  "'
Set-TextValue "D11" '="20.45"'
Set-TextValue "E11" '="  -1.58%  "'
Set-TextValue "D12" '="1.851.73"'
Set-TextValue "E12" '="  +1.43%  "'
Set-TextValue "D13" '="5.351"'
Set-TextValue "E13" '="  -1.21%  "'
Set-TextValue "D14" '="6.509"'
Set-TextValue "E14" '="  -2.68%  "'
Set-TextValue "D15" '="0.07051"'
Set-TextValue "E15" '="  -0.50%  "'
Set-TextValue "D16" '="91.20"'
Set-TextValue "E16" '="  -1.91%  "'
Set-TextValue "D17" '="1.003"'
Set-TextValue "E17" '="  +0.14%  "'
Set-TextValue "D18" '="0.000008703"'
Set-TextValue "E18" '="  -1.05%  "'
Set-TextValue "E19" '="  +0.16%  "'
Set-TextValue "D20" '="14.68"'
Set-TextValue "E20" '="  -2.27%  "'
Set-TextValue "D21" '="26.896.47"'
Set-TextValue "E21" '="  -1.73%  "'
Set-TextValue "D22" '="5.321"'
Set-TextValue "E22" '="  +0.04%  "'
Set-TextValue "D23" '="10.51"'
Set-TextValue "E23" '="  -4.31%  "'
Set-TextValue "D24" '="2.080.22"'
Set-TextValue "E24" '="  +1.35%  "'
Set-TextValue "D25" '="1.906"'
Set-TextValue "E25" '="  -1.73%  "'
Set-TextValue "E26" '="  +0.28%  "'
Set-TextValue "D27" '="18.36"'
Set-TextValue "E27" '="  -1.65%  "'
Set-TextValue "D28" '="2.138"'
Set-TextValue "E28" '="  -5.04%  "'
Set-TextValue "D29" '="5.302"'
Set-TextValue "E29" '="  -0.66%  "'
Set-TextValue "D30" '="115.81"'
Set-TextValue "E30" '="  -1.12%  "'
Set-TextValue "D31" '="0.08892"'
Set-TextValue "E31" '="  -0.71%  "'
Set-TextValue "D32" '="0.7516"'
Set-TextValue "E32" '="  -4.56%  "'
Set-TextValue "D33" '="1.151"'
Set-TextValue "E33" '="  -3.48%  "'
Set-TextValue "B34" '="HuobiToken"'
Set-TextValue "C34" '="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"'
Set-TextValue "D34" '="2.919"'
Set-TextValue "E34" '="  -0.70%  "'
Set-TextValue "B35" '="Filecoin"'
Set-TextValue "C35" '="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"'
Set-TextValue "D35" '="4.456"'
Set-TextValue "E35" '="  -1.78%  "'
Set-TextValue "E36" '="  +0.07%  "'
Set-TextValue "D37" '="1.097"'
Set-TextValue "E37" '="  -0.43%  "'
Set-TextValue "D38" '="0.01958"'
Set-TextValue "E38" '="  -0.89%  "'
Set-TextValue "D39" '="0.05251"'
Set-TextValue "E39" '="  +0.09%  "'
Set-TextValue "D40" '="2.420"'
Set-TextValue "E40" '="  +2.91%  "'
Set-TextValue "D41" '="2.928"'
Set-TextValue "D42" '="0.5303"'
Set-TextValue "E42" '="  -1.19%  "'
Set-TextValue "D43" '="7.161"'
Set-TextValue "E43" '="  -1.82%  "'
Set-TextValue "D44" '="0.1662"'
Set-TextValue "E44" '="  -2.47%  "'
Set-TextValue "D45" '="8.432"'
Set-TextValue "E45" '="  -2.58%  "'
Set-TextValue "D46" '="0.4932"'
Set-TextValue "E46" '="  -3.23%  "'
Set-TextValue "E47" '="  -3.34%  "'
Set-TextValue "E48" '="  +0.11%  "'
Set-TextValue "D49" '="1.668"'
Set-TextValue "E49" '="  -0.80%  "'
Set-TextValue "D50" '="102.99"'
Set-TextValue "E50" '="  -2.77%  "'
Set-TextValue "D51" '="0.06278"'
Set-TextValue "E51" '="  -1.67%  "'
